# edit.ps1 - apply commit "Add error codes Fix data storage issues Improve tests"
# to tests/data/test_data_01.xlsx
#
# Summary of changes:
#  - "train" sheet: populate rows 2-16 with NLP extraction results (idx/text/
#    subtext/span/entity/tag) for the 6 sample sentences.
#  - "config" sheet: drop the obsolete "nlp_name" row (row 3) and make the
#    "prepare_enabled" value an explicit boolean.
#  - Active sheet moves from "train" to "source".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# train sheet: append NLP-extracted entity rows
# ---------------------------------------------------------------------
$train = $wb.Worksheets.Item("train")

# row -> idx, text (col B)
$sentenceRows = @(
  @{ r = 2;  idx = 1;   text = "Uber blew through `$1 million a week" },
  @{ r = 5;  idx = 2;   text = "Android Pay expands to Canada" },
  @{ r = 8;  idx = 3;   text = "Spotify steps up Asia expansion" },
  @{ r = 11; idx = 4;   text = "Google Maps launches location sharing" },
  @{ r = 13; idx = 5;   text = "Google rebrands its business apps" },
  @{ r = 15; idx = 6;   text = "look what i found on google! 😂" }
)
foreach ($row in $sentenceRows) {
  $train.Cells.Item($row.r, 1).Value = $row.idx
  $train.Cells.Item($row.r, 2).Value = $row.text
}

# row -> idx, subtext (col C), span (col D, only row 3 has one), entity tag (col E)
$entityRows = @(
  @{ r = 3;  idx = 1.1; subtext = "Uber";        span = "0,4"; entity = "ORG" },
  @{ r = 4;  idx = 1.2; subtext = "`$1 million";  span = $null; entity = "MONEY" },
  @{ r = 6;  idx = 2.1; subtext = "Android Pay"; span = $null; entity = "PRODUCT" },
  @{ r = 7;  idx = 2.2; subtext = "Canada";      span = $null; entity = "GPE" },
  @{ r = 9;  idx = 3.1; subtext = "Spotify";     span = $null; entity = "ORG" },
  @{ r = 10; idx = 3.2; subtext = "Asia";        span = $null; entity = "LOC" },
  @{ r = 12; idx = 4.1; subtext = "Google Maps"; span = $null; entity = "PRODUCT" },
  @{ r = 14; idx = 5.1; subtext = "Google";      span = $null; entity = "ORG" },
  @{ r = 16; idx = 6.1; subtext = "google";      span = $null; entity = "PRODUCT" }
)
foreach ($row in $entityRows) {
  $train.Cells.Item($row.r, 1).Value = $row.idx
  $train.Cells.Item($row.r, 3).Value = $row.subtext
  if ($row.span -ne $null) {
    $train.Cells.Item($row.r, 4).Value = $row.span
  }
  $train.Cells.Item($row.r, 5).Value = $row.entity
}

# ---------------------------------------------------------------------
# config sheet: remove the obsolete "nlp_name" row, tidy prepare_enabled
# ---------------------------------------------------------------------
$config = $wb.Worksheets.Item("config")
$config.Rows.Item(3).Delete()

# After the delete, "prepare_enabled" is on row 4 -> store as explicit boolean
$config.Range("B4").Value = $true

# reset config sheet's selection back to A1 (its B1:B8-sized selection is
# stale now that the sheet only has 7 rows)
$null = $config.Range("A1").Select()

# ---------------------------------------------------------------------
# switch the active sheet from "train" to "source"
# ---------------------------------------------------------------------
$source = $wb.Worksheets.Item("source")
$source.Activate()
